$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New month column header - stored as text (like the other month headers
# B1:M1), not as a number. Format as Text first so the numeric-looking
# string "201912" is kept as text, then clear the number format again so
# the cell is left with the default style (matching its neighbours).
$ws.Range("N1").NumberFormat = "@"
$ws.Range("N1").Value = "201912"
$ws.Range("N1").ClearFormats()

# New month data values per row (row -> value); rows 7 and 26 have no data for this month
$values = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 0
    6  = 0
    8  = 1
    9  = 4
    10 = 0
    11 = 2
    12 = 2
    13 = 1
    14 = 1
    15 = 1
    16 = 3
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 1
    22 = 1
    23 = 0
    24 = 0
    25 = 0
    27 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("N$row").Value = $values[$row]
}
